# Apply the "include error count" update to the counters sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Zero-pad the single-digit months inside the embedded regex-mismatch
#    error messages (dates were re-run through the validator with a
#    leading zero on the month).
$ws.Range("C24").Value = "Error: Value '02/18/2011' does not match regex '^(((0?[1-9]|[12]\d|3[01])[\.\-\/](0?[13578]|1[02])[\.\-\/]((1[6-9]|[2-9]\d)?\d{2}|\d))|((0?[1-9]|[12]\d|30)[\.\-\/](0?[13456789]|1[012])[\.\-\/]((1[6-9]|[2-9]\d)?\d{2}|\d))|((0?[1-9]|1\d|2[0-8])[\.\-\/]0?2[\.\-\/]((1[6-9]|[2-9]\d)?\d{2}|\d))|(29[\.\-\/]0?2[\.\-\/]((1[6-9]|[2-9]\d)?(0[48]|[2468][048]|[13579][26])|((16|[2468][048]|[3579][26])00)|00|[048])))$'"
$ws.Range("C25").Value = "Error: Value '07/26/2010' does not match regex '^(((0?[1-9]|[12]\d|3[01])[\.\-\/](0?[13578]|1[02])[\.\-\/]((1[6-9]|[2-9]\d)?\d{2}|\d))|((0?[1-9]|[12]\d|30)[\.\-\/](0?[13456789]|1[012])[\.\-\/]((1[6-9]|[2-9]\d)?\d{2}|\d))|((0?[1-9]|1\d|2[0-8])[\.\-\/]0?2[\.\-\/]((1[6-9]|[2-9]\d)?\d{2}|\d))|(29[\.\-\/]0?2[\.\-\/]((1[6-9]|[2-9]\d)?(0[48]|[2468][048]|[13579][26])|((16|[2468][048]|[3579][26])00)|00|[048])))$'"
$ws.Range("C31").Value = "Error: Value '04/13/2012' does not match regex '^(((0?[1-9]|[12]\d|3[01])[\.\-\/](0?[13578]|1[02])[\.\-\/]((1[6-9]|[2-9]\d)?\d{2}|\d))|((0?[1-9]|[12]\d|30)[\.\-\/](0?[13456789]|1[012])[\.\-\/]((1[6-9]|[2-9]\d)?\d{2}|\d))|((0?[1-9]|1\d|2[0-8])[\.\-\/]0?2[\.\-\/]((1[6-9]|[2-9]\d)?\d{2}|\d))|(29[\.\-\/]0?2[\.\-\/]((1[6-9]|[2-9]\d)?(0[48]|[2468][048]|[13579][26])|((16|[2468][048]|[3579][26])00)|00|[048])))$'"
$ws.Range("C43").Value = "Error: Value '03/28/2012' does not match regex '^(((0?[1-9]|[12]\d|3[01])[\.\-\/](0?[13578]|1[02])[\.\-\/]((1[6-9]|[2-9]\d)?\d{2}|\d))|((0?[1-9]|[12]\d|30)[\.\-\/](0?[13456789]|1[012])[\.\-\/]((1[6-9]|[2-9]\d)?\d{2}|\d))|((0?[1-9]|1\d|2[0-8])[\.\-\/]0?2[\.\-\/]((1[6-9]|[2-9]\d)?\d{2}|\d))|(29[\.\-\/]0?2[\.\-\/]((1[6-9]|[2-9]\d)?(0[48]|[2468][048]|[13579][26])|((16|[2468][048]|[3579][26])00)|00|[048])))$'"
$ws.Range("C50").Value = "Error: Value '07/18/2012' does not match regex '^(((0?[1-9]|[12]\d|3[01])[\.\-\/](0?[13578]|1[02])[\.\-\/]((1[6-9]|[2-9]\d)?\d{2}|\d))|((0?[1-9]|[12]\d|30)[\.\-\/](0?[13456789]|1[012])[\.\-\/]((1[6-9]|[2-9]\d)?\d{2}|\d))|((0?[1-9]|1\d|2[0-8])[\.\-\/]0?2[\.\-\/]((1[6-9]|[2-9]\d)?\d{2}|\d))|(29[\.\-\/]0?2[\.\-\/]((1[6-9]|[2-9]\d)?(0[48]|[2468][048]|[13579][26])|((16|[2468][048]|[3579][26])00)|00|[048])))$'"

# 2. Four more "post" error rows showed up in the re-run, so insert four
#    blank rows right after the existing post/Completeness block (row 73)
#    -- this pushes the neighborhood (74-77 -> 78-81) and Location 1
#    (78-81 -> 82-85) blocks down by four rows, matching dimension C85.
$ws.Rows("74:77").Insert()

# 3. Populate the newly-inserted rows with the same "post" / mandatory
#    field completeness error as the existing post rows.
for ($r = 74; $r -le 77; $r++) {
    $ws.Cells.Item($r, 1).Value = "post"
    $ws.Cells.Item($r, 2).Value = "Completeness of Mandatory fields"
    $ws.Cells.Item($r, 3).Value = "Error: Mandatory field is BLANK or NULL. A value is required."
}

# 4. Every other "post" row (71, 73, 75, 77) now also reports a data-type
#    compliance error instead of the completeness one.
foreach ($r in 71, 73, 75, 77) {
    $ws.Cells.Item($r, 2).Value = "Meta Compliance (data type)"
    $ws.Cells.Item($r, 3).Value = "Error: Value '' is not an int. An int was expected"
}
